$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (strikeouts) column values replacing the old "Strike#" derived
# figures, recalculated from the pitch-by-pitch source data.
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 0
    24 = 3
    25 = 0
    26 = 0
    27 = 1
    28 = 2
    29 = 3
    30 = 3
    31 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
